$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.268.43"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.748.31"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "115.22"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "332.70"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.530"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +3.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.39"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.17"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.176.92"
$ws.Range("E15").Value = "  +3.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.717.08"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.887"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.156.87"
$ws.Range("E18").Value = "  +2.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.76"
$ws.Range("E19").Value = "  +4.56%  "
$ws.Range("E20").Value = "  +4.01%  "
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.55"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.09"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("E25").Value = "  +2.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.97"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.34"
$ws.Range("E28").Value = "  +2.32%  "
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.76"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.00"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.63"
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0826"
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.46"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.05"
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.11"
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "129.42"
$ws.Range("E40").Value = "  +3.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.75"
$ws.Range("E41").Value = "  +4.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0351"
$ws.Range("E42").Value = "  +10.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.113"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.29"
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.114.18"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.40"
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("E47").Value = "  +9.64%  "
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.56"
$ws.Range("E49").Value = "  +3.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.06"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("E51").Value = "  +9.78%  "
